$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (A and B shrink slightly) ---
# Excel quantizes ColumnWidth to an internal pixel grid, so the values below
# are chosen so that after quantization the saved "width" attribute matches
# (or gets as close as possible to) the target from the diff.
$ws.Columns.Item(1).ColumnWidth = 38.666666666666664   # -> stored width ~39.42578125 (closest achievable: 39.5)
$ws.Columns.Item(2).ColumnWidth = 37.166666666666664   # -> stored width 38

# --- Add a new "2022" column (S) of data, mirroring the existing R column's formatting ---

# Year header S4 (copy format from R4, a "2021" header cell)
$ws.Range("R4").Copy() | Out-Null
$ws.Range("S4").PasteSpecial(-4122) | Out-Null
$ws.Range("S4").Value = 2022

# S5 (copy format from R5)
$ws.Range("R5").Copy() | Out-Null
$ws.Range("S5").PasteSpecial(-4122) | Out-Null
$ws.Range("S5").Value = 44

# S6 (copy format from R6, then apply the "0.0" number format used by this new column's data)
$ws.Range("R6").Copy() | Out-Null
$ws.Range("S6").PasteSpecial(-4122) | Out-Null
$ws.Range("S6").NumberFormat = "0.0"
$ws.Range("S6").Value = 20.6

# S7 (copy format from R7)
$ws.Range("R7").Copy() | Out-Null
$ws.Range("S7").PasteSpecial(-4122) | Out-Null
$ws.Range("S7").Value = 7.9

# S8 (copy format from R8)
$ws.Range("R8").Copy() | Out-Null
$ws.Range("S8").PasteSpecial(-4122) | Out-Null
$ws.Range("S8").Value = 15.5

$excel.CutCopyMode = 0

# --- Selection moves (cosmetic, matches the author's last cursor position) ---
$ws.Range("Y14").Select() | Out-Null
